# Updates Price (D) and Volume(1h) (E) figures for several coins, as
# published by the "Updated symbol list" GitHub Actions job.
# Values are written with a leading apostrophe so Excel keeps them as
# literal text (matching the original inline-string cells) instead of
# auto-converting numeric-/percent-looking text into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.37"
$ws.Range("E2").Value = "'0.12%"
$ws.Range("D3").Value = "'26.75"
$ws.Range("E3").Value = "'-2.20%"
$ws.Range("D4").Value = "'4.701"
$ws.Range("E4").Value = "'-0.03%"
$ws.Range("D5").Value = "'0.06193"
$ws.Range("E5").Value = "'1.69%"
$ws.Range("D6").Value = "'6.747"
$ws.Range("E6").Value = "'1.08%"
$ws.Range("D7").Value = "'0.8494"
$ws.Range("E7").Value = "'0.35%"
$ws.Range("D8").Value = "'0.9120"
$ws.Range("E8").Value = "'-1.63%"
$ws.Range("D9").Value = "'0.1401"
$ws.Range("E9").Value = "'-0.21%"
$ws.Range("D10").Value = "'0.05083"
$ws.Range("E10").Value = "'5.81%"
$ws.Range("D11").Value = "'0.07101"
$ws.Range("E11").Value = "'0.06%"
$ws.Range("D12").Value = "'0.03107"
$ws.Range("D13").Value = "'0.09041"
$ws.Range("E13").Value = "'-0.24%"
$ws.Range("E14").Value = "'-0.12%"
$ws.Range("D15").Value = "'0.0006163"
$ws.Range("E15").Value = "'1.68%"
$ws.Range("D16").Value = "'0.005944"
$ws.Range("E16").Value = "'-2.77%"
$ws.Range("D17").Value = "'3.447"
$ws.Range("E17").Value = "'-0.08%"
$ws.Range("D18").Value = "'3.173"
$ws.Range("E18").Value = "'0.78%"
$ws.Range("D20").Value = "'0.3072"
$ws.Range("E20").Value = "'-1.21%"
$ws.Range("E21").Value = "'0.40%"
$ws.Range("D22").Value = "'4.115"
$ws.Range("E22").Value = "'0.39%"
$ws.Range("D23").Value = "'0.04256"
$ws.Range("E23").Value = "'0.48%"
$ws.Range("D24").Value = "'0.001185"
$ws.Range("E24").Value = "'-3.08%"
$ws.Range("E25").Value = "'6.94%"
$ws.Range("D40").Value = "'0.03966"
$ws.Range("E40").Value = "'2.37%"
$ws.Range("D41").Value = "'0.1111"
$ws.Range("E41").Value = "'-0.20%"
$ws.Range("D42").Value = "'0.004140"
$ws.Range("E42").Value = "'1.53%"
$ws.Range("D43").Value = "'0.002143"
$ws.Range("E43").Value = "'-3.33%"
$ws.Range("D44").Value = "'0.01325"
$ws.Range("E44").Value = "'-18.41%"
$ws.Range("D45").Value = "'0.00005164"
$ws.Range("E45").Value = "'0.31%"
$ws.Range("D48").Value = "'0.2481"
$ws.Range("E48").Value = "'80.98%"
